# Daily data-refresh update: append the 2020-06-18 row to the
# "Condicion_Pacientes" table on Hoja1 (table + sheet both grow from
# A1:F97 to A1:F98).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data row (97) down into the
# new row (98) so the new cells pick up the same date / centered number
# styles used throughout the table, then overwrite with the new values.
$ws.Range("A97:F97").Copy()
$ws.Range("A98:F98").PasteSpecial(-4122)

$ws.Range("A98").Value = 44000
$ws.Range("B98").Value = 1111
$ws.Range("C98").Value = 440
$ws.Range("D98").Value = 590
$ws.Range("E98").Value = 295
$ws.Range("F98").Value = 48

# Grow the table (and its autoFilter range) to include the new row.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.Resize($ws.Range("A1:F98"))

# Match the author's final selection on the newly-added last cell.
$ws.Range("F98").Select() | Out-Null
